$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "within the ML4 course on " and "Coursera" were two separate
# runs (the second wrapped in proofErr spell-check markers). Merge them
# into a single run/text and drop the proofErr markers.
# ---------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Type of event or data within the ML4 course on Coursera", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$cellPara = $findRange.Paragraphs.Item(1)
$paraStart = $cellPara.Range.Start
$paraEnd = $cellPara.Range.End

# Extend by one so the replaced range swallows the paragraph's own end-of
# -paragraph mark too -- InsertXML then substitutes the whole <w:p> in
# place instead of inserting a sibling paragraph next to it.
$wholePara = $d.Range($paraStart, $paraEnd + 1)

$mergedParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="74CC17A1" w14:textId="2FB294C3" w:rsidR="002C4154" w:rsidRPr="003C3BC6" w:rsidRDefault="002422A8" w:rsidP="00BB5891"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve">Type of event or data </w:t></w:r><w:r w:rsidR="00BB5891"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t>within the ML4 course on Coursera</w:t></w:r></w:p>'

$wholePara.InsertXML($mergedParaXml)

# ---------------------------------------------------------------------
# Change 2: drop an extra empty paragraph right after the last table, to
# improve the placement of the figures that follow it.
# ---------------------------------------------------------------------
$lastTable = $d.Tables.Item($d.Tables.Count)
$afterTable = $d.Range($lastTable.Range.End, $lastTable.Range.End)
$afterTable.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Output "ok"
